$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999139175599,
    0.99542249440631603,
    0.97677190956773585,
    0.96484215510203297,
    0.95325842359145574,
    0.92681712393285021,
    0.92084082190625005,
    0.91278545097208785,
    0.90281728866262057,
    0.89376591300209163,
    0.89245966678232214,
    0.8902307661794614,
    0.88147178116315605,
    0.87839144100135691,
    0.87684062450847122,
    0.87557864220962278,
    0.87444700500906847,
    0.87457514460937646,
    0.99380917045483752,
    0.98145316569247187,
    0.98005469764754749,
    0.97879019465230632,
    0.97508272059941903,
    0.9730818113259615,
    0.96662492940249689,
    0.92424180827153746,
    0.92090014273189591,
    0.90827866721642669,
    0.89955994522625837,
    0.89546187133212951,
    0.88780838627239089,
    0.88612910499261377,
    0.88560910686182459
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
